$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the scrolled viewport (topLeftCell="A7") so the view resets to the top-left.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1

# --- Insert "numCommande" row right after "idCommande" (row 17) ---
$ws.Rows(18).Insert()
$ws.Range("A18").Value = "numCommande"
$ws.Range("B18").Value = "Numéro de la commande"
$ws.Range("C18").Value = "int"
$ws.Range("E18").Value = "N"

# --- Insert "numLigneCommande" row right after "idLigneCommande" (now row 21) ---
$ws.Rows(22).Insert()
$ws.Range("A22").Value = "numLigneCommande"
$ws.Range("B22").Value = "Numéro de la ligne de la commande"
$ws.Range("C22").Value = "int"
$ws.Range("E22").Value = "N"
